$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: new record - JESUS GABRIEL JIMENEZ TUÑON (moved up from row 19)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143357687"
$ws.Range("D16").Value = "JESUS GABRIEL JIMENEZ TUÑON"
$ws.Range("E16").Value = "1808"
$ws.Range("F16").Value = 21333
$ws.Range("G16").Value = 1000000

# Row 17: EDILMER JOS ARRIETA PINEDA - periodo 2002, updated salario
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047369991"
$ws.Range("D17").Value = "EDILMER JOS ARRIETA PINEDA"
$ws.Range("E17").Value = "2002"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 877803

# Row 18: EDILMER JOS ARRIETA PINEDA - periodo 2003, updated salario
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047369991"
$ws.Range("D18").Value = "EDILMER JOS ARRIETA PINEDA"
$ws.Range("E18").Value = "2003"
$ws.Range("F18").Value = 33125
$ws.Range("G18").Value = 877803

# Row 19: EDILMER JOS ARRIETA PINEDA - periodo 2004, updated salario
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047369991"
$ws.Range("D19").Value = "EDILMER JOS ARRIETA PINEDA"
$ws.Range("E19").Value = "2004"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 877803
